$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.064.25"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.058.44"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.10"
$ws.Range("E5").Value = "  -1.65%  "
$ws.Range("E6").Value = "  -1.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.10"
$ws.Range("E7").Value = "  -1.51%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.380"
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0776"
$ws.Range("E10").Value = "  -2.27%  "
$ws.Range("E11").Value = "  +2.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.53"
$ws.Range("E12").Value = "  -3.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.891"
$ws.Range("E13").Value = "  +7.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.360.03"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.74"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.052.89"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.22"
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.029.94"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.24"
$ws.Range("E19").Value = "  -1.81%  "
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.49"
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "239.39"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  +1.61%  "
$ws.Range("E25").Value = "  +4.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.38"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  -4.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.19"
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.45"
$ws.Range("E29").Value = "  +12.54%  "
$ws.Range("E30").Value = "  -1.36%  "
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.71"
$ws.Range("E32").Value = "  +3.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0618"
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("E34").Value = "  +3.20%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.83"
$ws.Range("E36").Value = "  +5.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0840"
$ws.Range("E37").Value = "  -6.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.34"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.28"
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("E42").Value = "  +2.02%  "
$ws.Range("E43").Value = "  -10.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.14"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.09"
$ws.Range("E45").Value = "  -4.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.305.81"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.38"
$ws.Range("E47").Value = "  -4.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.87"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.84"
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.248.09"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.72"
$ws.Range("E51").Value = "  +2.05%  "
